$wb = $excel.ActiveWorkbook

$openQuote = [char]0x201C
$closeQuote = [char]0x201D

for ($i = 1; $i -le 20; $i++) {
    $sheetIndex = $i + 3
    $ws = $wb.Worksheets.Item($sheetIndex)
    $oldName = "Consequence of " + $openQuote + "A" + $i + $closeQuote
    $newName = $openQuote + "A" + $i + $closeQuote + " consequence"
    if ($ws.Name -eq $oldName) {
        $ws.Name = $newName
    }
}
